$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 21.25004918957847
$ws.Cells.Item(2, 3).Value = 2.894004852692718
$ws.Cells.Item(2, 4).Value = 28.71952398889674
$ws.Cells.Item(2, 5).Value = 64.72726527080729
$ws.Cells.Item(2, 6).Value = 27.45237635326557
$ws.Cells.Item(2, 7).Value = 0.8035723718204869
$ws.Cells.Item(2, 8).Value = 15.55334958792415
$ws.Cells.Item(2, 9).Value = 11.76492207051695
$ws.Cells.Item(2, 10).Value = 15.69700588368856

$ws.Cells.Item(3, 2).Value = 12.91032764359433
$ws.Cells.Item(3, 3).Value = 7.180446301491812
$ws.Cells.Item(3, 4).Value = 31.25189154971744
$ws.Cells.Item(3, 5).Value = 64.87660640106319
$ws.Cells.Item(3, 6).Value = 25.65096239551387
$ws.Cells.Item(3, 7).Value = 0.803803684545715
$ws.Cells.Item(3, 8).Value = 15.20171594488467
$ws.Cells.Item(3, 9).Value = 11.91387822353283
$ws.Cells.Item(3, 10).Value = 15.22430741713245

$ws.Cells.Item(4, 2).Value = 9.679497218986368
$ws.Cells.Item(4, 3).Value = 10.35036420107674
$ws.Cells.Item(4, 4).Value = 33.01111778196978
$ws.Cells.Item(4, 5).Value = 64.71543602826131
$ws.Cells.Item(4, 6).Value = 24.2177523677013
$ws.Cells.Item(4, 7).Value = 0.8039597838511461
$ws.Cells.Item(4, 8).Value = 14.96464539490871
$ws.Cells.Item(4, 9).Value = 11.74775628332204
$ws.Cells.Item(4, 10).Value = 14.83312558573151

$ws.Cells.Item(5, 2).Value = 9.97014183677604
$ws.Cells.Item(5, 3).Value = 11.71575512395725
$ws.Cells.Item(5, 4).Value = 33.62026150034059
$ws.Cells.Item(5, 5).Value = 64.28998674927752
$ws.Cells.Item(5, 6).Value = 23.45870343655098
$ws.Cells.Item(5, 7).Value = 0.8040525171553548
$ws.Cells.Item(5, 8).Value = 14.82361822660791
$ws.Cells.Item(5, 9).Value = 11.31390889227738
$ws.Cells.Item(5, 10).Value = 14.57460341736893

$ws.Cells.Item(6, 2).Value = 10.3584346003494
$ws.Cells.Item(6, 3).Value = 12.19066688394368
$ws.Cells.Item(6, 4).Value = 33.76274696693466
$ws.Cells.Item(6, 5).Value = 63.97566787628074
$ws.Cells.Item(6, 6).Value = 23.12776794563901
$ws.Cells.Item(6, 7).Value = 0.8041002311544542
$ws.Cells.Item(6, 8).Value = 14.75108752742528
$ws.Cells.Item(6, 9).Value = 10.99373759213813
$ws.Cells.Item(6, 10).Value = 14.44009105561576

$ws.Cells.Item(7, 2).Value = 10.3584346003494
$ws.Cells.Item(7, 3).Value = 12.19066688394368
$ws.Cells.Item(7, 4).Value = 33.76274696693466
$ws.Cells.Item(7, 5).Value = 63.97566787628074
$ws.Cells.Item(7, 6).Value = 23.12776794563901
$ws.Cells.Item(7, 7).Value = 0.8041002311544542
$ws.Cells.Item(7, 8).Value = 14.75108752742528
$ws.Cells.Item(7, 9).Value = 10.99373759213813
$ws.Cells.Item(7, 10).Value = 14.44009105561576

$ws.Cells.Item(8, 2).Value = 10.14873937099412
$ws.Cells.Item(8, 3).Value = 11.94339324343611
$ws.Cells.Item(8, 4).Value = 33.59598109348281
$ws.Cells.Item(8, 5).Value = 63.91687620202312
$ws.Cells.Item(8, 6).Value = 23.20931105439569
$ws.Cells.Item(8, 7).Value = 0.804098633370848
$ws.Cells.Item(8, 8).Value = 14.7536587447172
$ws.Cells.Item(8, 9).Value = 10.93407342304825
$ws.Cells.Item(8, 10).Value = 14.44800621418932

$ws.Cells.Item(9, 2).Value = 10.18153019857932
$ws.Cells.Item(9, 3).Value = 9.105897723562899
$ws.Cells.Item(9, 4).Value = 31.92158654744843
$ws.Cells.Item(9, 5).Value = 63.82480882473538
$ws.Cells.Item(9, 6).Value = 24.38890163713014
$ws.Cells.Item(9, 7).Value = 0.8040034960364142
$ws.Cells.Item(9, 8).Value = 14.8982783944484
$ws.Cells.Item(9, 9).Value = 10.84232482539293
$ws.Cells.Item(9, 10).Value = 14.71440285107261

$ws.Cells.Item(10, 2).Value = 21.60908680432514
$ws.Cells.Item(10, 3).Value = 2.620652844588401
$ws.Cells.Item(10, 4).Value = 28.30620006851873
$ws.Cells.Item(10, 5).Value = 64.11590622003195
$ws.Cells.Item(10, 6).Value = 27.31823870787008
$ws.Cells.Item(10, 7).Value = 0.8036206907430595
$ws.Cells.Item(10, 8).Value = 15.48017358462991
$ws.Cells.Item(10, 9).Value = 11.14317749104015
$ws.Cells.Item(10, 10).Value = 15.55592419136711

$ws.Cells.Item(11, 2).Value = 35.96340458994777
$ws.Cells.Item(11, 3).Value = 3.998623159941701
$ws.Cells.Item(11, 4).Value = 24.80729071909315
$ws.Cells.Item(11, 5).Value = 64.31782551968915
$ws.Cells.Item(11, 6).Value = 30.13175860908996
$ws.Cells.Item(11, 7).Value = 0.8031116548094855
$ws.Cells.Item(11, 8).Value = 16.2543198403162
$ws.Cells.Item(11, 9).Value = 11.35321411401343
$ws.Cells.Item(11, 10).Value = 16.46417785932586

$ws.Cells.Item(12, 2).Value = 42.78201381643699
$ws.Cells.Item(12, 3).Value = 6.844962290645418
$ws.Cells.Item(12, 4).Value = 22.9723535430511
$ws.Cells.Item(12, 5).Value = 64.2248662800219
$ws.Cells.Item(12, 6).Value = 31.4729681098063
$ws.Cells.Item(12, 7).Value = 0.8028294841862911
$ws.Cells.Item(12, 8).Value = 16.63885378937666
$ws.Cells.Item(12, 9).Value = 11.26100131462553
$ws.Cells.Item(12, 10).Value = 16.90454976151511

$ws.Cells.Item(13, 2).Value = 44.71484321558669
$ws.Cells.Item(13, 3).Value = 7.331762562171729
$ws.Cells.Item(13, 4).Value = 21.525746664147
$ws.Cells.Item(13, 5).Value = 64.12208190309023
$ws.Cells.Item(13, 6).Value = 32.51298011629294
$ws.Cells.Item(13, 7).Value = 0.8026383913957232
$ws.Cells.Item(13, 8).Value = 16.76152620307452
$ws.Cells.Item(13, 9).Value = 11.15818699286314
$ws.Cells.Item(13, 10).Value = 17.23596960831005

$ws.Cells.Item(14, 2).Value = 42.51013335544093
$ws.Cells.Item(14, 3).Value = 6.117021462654916
$ws.Cells.Item(14, 4).Value = 21.15081052466358
$ws.Cells.Item(14, 5).Value = 63.97626703671029
$ws.Cells.Item(14, 6).Value = 32.69357590160311
$ws.Cells.Item(14, 7).Value = 0.802639594981338
$ws.Cells.Item(14, 8).Value = 16.80430327510293
$ws.Cells.Item(14, 9).Value = 11.01024835274076
$ws.Cells.Item(14, 10).Value = 17.26042288175982

$ws.Cells.Item(15, 2).Value = 40.13138671346117
$ws.Cells.Item(15, 3).Value = 5.026749999261418
$ws.Cells.Item(15, 4).Value = 21.28302630779413
$ws.Cells.Item(15, 5).Value = 63.76046960007713
$ws.Cells.Item(15, 6).Value = 32.42972522919675
$ws.Cells.Item(15, 7).Value = 0.8027169742672819
$ws.Cells.Item(15, 8).Value = 16.75994288772146
$ws.Cells.Item(15, 9).Value = 10.79053987928663
$ws.Cells.Item(15, 10).Value = 17.12857295780628

$ws.Cells.Item(16, 2).Value = 39.24629953925729
$ws.Cells.Item(16, 3).Value = 4.571881288312859
$ws.Cells.Item(16, 4).Value = 21.4837303893533
$ws.Cells.Item(16, 5).Value = 63.9718326719511
$ws.Cells.Item(16, 6).Value = 32.43222462941379
$ws.Cells.Item(16, 7).Value = 0.8027166936180201
$ws.Cells.Item(16, 8).Value = 16.79362326639581
$ws.Cells.Item(16, 9).Value = 11.00527636592407
$ws.Cells.Item(16, 10).Value = 17.15839448815192

$ws.Cells.Item(17, 2).Value = 38.4475885612055
$ws.Cells.Item(17, 3).Value = 4.15484224715394
$ws.Cells.Item(17, 4).Value = 21.72763845282828
$ws.Cells.Item(17, 5).Value = 64.24621479399441
$ws.Cells.Item(17, 6).Value = 32.44837863882517
$ws.Cells.Item(17, 7).Value = 0.8027103938138469
$ws.Cells.Item(17, 8).Value = 16.83012774910791
$ws.Cells.Item(17, 9).Value = 11.2840780677501
$ws.Cells.Item(17, 10).Value = 17.20328707664975

$ws.Cells.Item(18, 2).Value = 34.30676021674427
$ws.Cells.Item(18, 3).Value = 2.461203395444084
$ws.Cells.Item(18, 4).Value = 22.80651335650802
$ws.Cells.Item(18, 5).Value = 64.22106499158325
$ws.Cells.Item(18, 6).Value = 31.59667313744858
$ws.Cells.Item(18, 7).Value = 0.8028864103942324
$ws.Cells.Item(18, 8).Value = 16.5501393911906
$ws.Cells.Item(18, 9).Value = 11.25715518253381
$ws.Cells.Item(18, 10).Value = 16.90741996719354

$ws.Cells.Item(19, 2).Value = 32.47660520622198
$ws.Cells.Item(19, 3).Value = 1.805748162548659
$ws.Cells.Item(19, 4).Value = 23.61016235901085
$ws.Cells.Item(19, 5).Value = 64.34466919717065
$ws.Cells.Item(19, 6).Value = 31.0699473445499
$ws.Cells.Item(19, 7).Value = 0.8029807701998505
$ws.Cells.Item(19, 8).Value = 16.33611874443728
$ws.Cells.Item(19, 9).Value = 11.38191899645494
$ws.Cells.Item(19, 10).Value = 16.7527961260207

$ws.Cells.Item(20, 2).Value = 36.8314265070984
$ws.Cells.Item(20, 3).Value = 3.988003953330698
$ws.Cells.Item(20, 4).Value = 23.58074092501053
$ws.Cells.Item(20, 5).Value = 64.32050420823153
$ws.Cells.Item(20, 6).Value = 31.07502308531781
$ws.Cells.Item(20, 7).Value = 0.8029521284003283
$ws.Cells.Item(20, 8).Value = 16.29193485654722
$ws.Cells.Item(20, 9).Value = 11.35743542684315
$ws.Cells.Item(20, 10).Value = 16.76407338034473

$ws.Cells.Item(21, 2).Value = 47.89935755655275
$ws.Cells.Item(21, 3).Value = 8.985080613704515
$ws.Cells.Item(21, 4).Value = 21.92480719668604
$ws.Cells.Item(21, 5).Value = 64.45701456332047
$ws.Cells.Item(21, 6).Value = 32.45529132187305
$ws.Cells.Item(21, 7).Value = 0.8025900455085839
$ws.Cells.Item(21, 8).Value = 17.0201041155886
$ws.Cells.Item(21, 9).Value = 11.49851977318397
$ws.Cells.Item(21, 10).Value = 17.2907316068616

$ws.Cells.Item(22, 2).Value = 54.94886710915529
$ws.Cells.Item(22, 3).Value = 11.92641056332669
$ws.Cells.Item(22, 4).Value = 20.49367735581331
$ws.Cells.Item(22, 5).Value = 64.76816533872685
$ws.Cells.Item(22, 6).Value = 33.79970625650883
$ws.Cells.Item(22, 7).Value = 0.8022355275288584
$ws.Cells.Item(22, 8).Value = 17.58469434207144
$ws.Cells.Item(22, 9).Value = 11.81679440829343
$ws.Cells.Item(22, 10).Value = 17.83895110264013

$ws.Cells.Item(23, 2).Value = 59.93827501442647
$ws.Cells.Item(23, 3).Value = 13.96319389979013
$ws.Cells.Item(23, 4).Value = 19.49502836351401
$ws.Cells.Item(23, 5).Value = 65.10457009485863
$ws.Cells.Item(23, 6).Value = 34.83248971952859
$ws.Cells.Item(23, 7).Value = 0.8019482501542721
$ws.Cells.Item(23, 8).Value = 18.02180875813952
$ws.Cells.Item(23, 9).Value = 12.16009977088045
$ws.Cells.Item(23, 10).Value = 18.28677549478493

$ws.Cells.Item(24, 2).Value = 61.58530534793537
$ws.Cells.Item(24, 3).Value = 14.61375400970072
$ws.Cells.Item(24, 4).Value = 19.2321711853004
$ws.Cells.Item(24, 5).Value = 65.34177410030624
$ws.Cells.Item(24, 6).Value = 35.21727320327172
$ws.Cells.Item(24, 7).Value = 0.8018379772456773
$ws.Cells.Item(24, 8).Value = 18.18981620445143
$ws.Cells.Item(24, 9).Value = 12.4016236398327
$ws.Cells.Item(24, 10).Value = 18.47194534792943

$ws.Cells.Item(25, 2).Value = 50.31541926115888
$ws.Cells.Item(25, 3).Value = 9.851889618845453
$ws.Cells.Item(25, 4).Value = 21.97773341698525
$ws.Cells.Item(25, 5).Value = 65.51972925621712
$ws.Cells.Item(25, 6).Value = 33.20888254968744
$ws.Cells.Item(25, 7).Value = 0.8023863592409626
$ws.Cells.Item(25, 8).Value = 17.35617218135903
$ws.Cells.Item(25, 9).Value = 12.57954909938501
$ws.Cells.Item(25, 10).Value = 17.72354820145952

Write-Output "applied loading_percent updates"